# frame work data setup
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the new row of test-data (TC_Name / Destination pair)
$ws.Range("A2").Value = "enterDetailsHomePage"
$ws.Range("B2").Value = "Mum"

# Column A needs to be a bit wider to show the longer value that was just added
$ws.Columns.Item(1).ColumnWidth = 31.666666666666668

# Move the active selection onto the newly populated column B
$ws.Range("B1").Select()

$wb.Save()
